# "Doing Updates for Financials"
# Update a handful of cells on the FPAFY sheet with refreshed financial figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Earnings Before Interest And Taxes (row 21) - 2011 figure (col J) is no longer available
$ws.Range("J21").Value = "NA"

# Depreciation (row 83) - 2011 figure (col J) is no longer available
$ws.Range("J83").Value = "NA"

# Capital Expenditures (row 91) - refreshed figures across all years
$ws.Range("D91").Value = -614100
$ws.Range("E91").Value = -319500
$ws.Range("F91").Value = -436800
$ws.Range("G91").Value = -482800
$ws.Range("H91").Value = -768000
$ws.Range("I91").Value = -541000
$ws.Range("J91").Value = -357300

# Total Cash Flows From Investing Activities (row 94)
$ws.Range("E94").Value = -593700
$ws.Range("J94").Value = "NA"

# Total Cash Flows From Financing Activities (row 100) - 2011 figure (col J) is no longer available
$ws.Range("J100").Value = "NA"

# Effect Of Exchange Rate Changes (row 101) - 2011 figure (col J) is no longer available
$ws.Range("J101").Value = "NA"
